# Apply the "use AbilityValues instead of AbilitySpecial" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: rename the "AbilitySpecial[{]" marker cell to "AbilityValues[{]"
$ws.Range("D2").Value = "AbilityValues[{]"

# Row 2: the var-name placeholder cells E2:M2 go from "01".."09" to "1".."9"
# (still stored as text thanks to the existing text number-format style)
$ws.Range("E2").Value = "1"
$ws.Range("F2").Value = "2"
$ws.Range("G2").Value = "3"
$ws.Range("H2").Value = "4"
$ws.Range("I2").Value = "5"
$ws.Range("J2").Value = "6"
$ws.Range("K2").Value = "7"
$ws.Range("L2").Value = "8"
$ws.Range("M2").Value = "9"

# Remove the now-unused instructional comment cell in row 1
$ws.Range("E1").ClearContents()

# Row 3 (the only remaining data row): add ScriptFile path, fix the
# "time1 1 2 3 4" typo, and clear the leftover numeric cells.
$ws.Range("C3").Value = "items/item_kv_generator_test.lua"
$ws.Range("F3").Value = "1 1 2 3 4"
$ws.Range("I3").ClearContents()
$ws.Range("J3").ClearContents()

# Row 4 was a duplicate test row ("item_kv_generator_test1") that got
# removed entirely in this revision.
$ws.Rows("4:4").Delete()

# Match the saved cursor position recorded in the workbook after the edit.
$ws.Range("D4").Select() | Out-Null
